# Auto-generated edit script: updates crafting leve profit/price data
# as refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2755.7778
$ws.Range("I86").Value = 2500.25
$ws.Range("J86").Value = 4800
$ws.Range("K86").Value = 2500.25
$ws.Range("L86").Value = 4800
$ws.Range("M86").Value = -1377.25
$ws.Range("N86").Value = -7046
$ws.Range("H88").Value = 2609.0715
$ws.Range("I88").Value = 1427
$ws.Range("J88").Value = 2931.4546
$ws.Range("K88").Value = 1427
$ws.Range("L88").Value = 2931.4546
$ws.Range("M88").Value = -1021
$ws.Range("N88").Value = -3743.4546
$ws.Range("H89").Value = 2755.7778
$ws.Range("I89").Value = 2500.25
$ws.Range("J89").Value = 4800
$ws.Range("K89").Value = 12501.25
$ws.Range("L89").Value = 24000
$ws.Range("M89").Value = -6885.25
$ws.Range("N89").Value = -35232
$ws.Range("H91").Value = 2609.0715
$ws.Range("I91").Value = 1427
$ws.Range("J91").Value = 2931.4546
$ws.Range("K91").Value = 1427
$ws.Range("L91").Value = 2931.4546
$ws.Range("M91").Value = -23
$ws.Range("N91").Value = -5739.4546
$ws.Range("H112").Value = 166668140
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 166668140
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 500004420
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -500006636
$ws.Range("H138").Value = 5459.6
$ws.Range("J138").Value = 6337.838
$ws.Range("L138").Value = 19013.514
$ws.Range("N138").Value = -29293.514

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11187.283
$ws.Range("I32").Value = 5289.2964
$ws.Range("J32").Value = 35686.617
$ws.Range("K32").Value = 5289.2964
$ws.Range("L32").Value = 35686.617
$ws.Range("M32").Value = -5002.2964
$ws.Range("N32").Value = -36260.617
$ws.Range("H74").Value = 6999.4443
$ws.Range("I74").Value = 899.3333
$ws.Range("J74").Value = 37500
$ws.Range("K74").Value = 899.3333
$ws.Range("L74").Value = 37500
$ws.Range("M74").Value = -25.33330000000001
$ws.Range("N74").Value = -39248
$ws.Range("H77").Value = 6999.4443
$ws.Range("I77").Value = 899.3333
$ws.Range("J77").Value = 37500
$ws.Range("K77").Value = 4496.6665
$ws.Range("L77").Value = 187500
$ws.Range("M77").Value = -128.6665000000003
$ws.Range("N77").Value = -196236
$ws.Range("H88").Value = 3971.1177
$ws.Range("I88").Value = 2401.5
$ws.Range("J88").Value = 5366.3335
$ws.Range("K88").Value = 2401.5
$ws.Range("L88").Value = 5366.3335
$ws.Range("M88").Value = -1995.5
$ws.Range("N88").Value = -6178.3335
$ws.Range("H91").Value = 3971.1177
$ws.Range("I91").Value = 2401.5
$ws.Range("J91").Value = 5366.3335
$ws.Range("K91").Value = 2401.5
$ws.Range("L91").Value = 5366.3335
$ws.Range("M91").Value = -997.5
$ws.Range("N91").Value = -8174.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 857.0606
$ws.Range("I20").Value = 594.5714
$ws.Range("J20").Value = 1316.4166
$ws.Range("K20").Value = 594.5714
$ws.Range("L20").Value = 1316.4166
$ws.Range("M20").Value = -347.5714
$ws.Range("N20").Value = -1810.4166
$ws.Range("H94").Value = 1078.4
$ws.Range("I94").Value = 1048.25
$ws.Range("J94").Value = 1138.7
$ws.Range("K94").Value = 1048.25
$ws.Range("L94").Value = 1138.7
$ws.Range("M94").Value = -597.25
$ws.Range("N94").Value = -2040.7
$ws.Range("H107").Value = 883026.8
$ws.Range("I107").Value = 1175410.8
$ws.Range("J107").Value = 5875
$ws.Range("K107").Value = 1175410.8
$ws.Range("L107").Value = 5875
$ws.Range("M107").Value = -1173490.8
$ws.Range("N107").Value = -9715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1263.2028
$ws.Range("I31").Value = 1080.4546
$ws.Range("J31").Value = 1410.2927
$ws.Range("K31").Value = 1080.4546
$ws.Range("L31").Value = 1410.2927
$ws.Range("M31").Value = -785.4546
$ws.Range("N31").Value = -2000.2927
$ws.Range("H34").Value = 1263.2028
$ws.Range("I34").Value = 1080.4546
$ws.Range("J34").Value = 1410.2927
$ws.Range("K34").Value = 1080.4546
$ws.Range("L34").Value = 1410.2927
$ws.Range("M34").Value = -878.4546
$ws.Range("N34").Value = -1814.2927

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1439.8
$ws.Range("I14").Value = 1439.8
$ws.Range("K14").Value = 4319.4
$ws.Range("M14").Value = -4146.4
$ws.Range("H49").Value = 2500
$ws.Range("J49").Value = 2500
$ws.Range("L49").Value = 7500
$ws.Range("N49").Value = -7812
$ws.Range("H57").Value = 2324.5
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 2766
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 8298
$ws.Range("M57").Value = -2441
$ws.Range("N57").Value = -9416
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 9000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -11122
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 27000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -37608
$ws.Range("H102").Value = 6004.3335
$ws.Range("I102").Value = 1013
$ws.Range("J102").Value = 8500
$ws.Range("K102").Value = 3039
$ws.Range("L102").Value = 25500
$ws.Range("M102").Value = -605
$ws.Range("N102").Value = -30368
$ws.Range("H107").Value = 59554.47
$ws.Range("J107").Value = 334834.66
$ws.Range("L107").Value = 1004503.98
$ws.Range("N107").Value = -1008343.98
$ws.Range("H110").Value = 2260
$ws.Range("I110").Value = 1750
$ws.Range("K110").Value = 5250
$ws.Range("M110").Value = -1160
$ws.Range("H131").Value = 1669504
$ws.Range("J131").Value = 2085521.9
$ws.Range("L131").Value = 6256565.699999999
$ws.Range("N131").Value = -6266645.699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1223.8572
$ws.Range("I97").Value = 1244.1072
$ws.Range("J97").Value = 1142.8572
$ws.Range("K97").Value = 1244.1072
$ws.Range("L97").Value = 1142.8572
$ws.Range("M97").Value = -748.1071999999999
$ws.Range("N97").Value = -2134.8572
$ws.Range("H132").Value = 1935.4783
$ws.Range("I132").Value = 1834.0952
$ws.Range("K132").Value = 5502.2856
$ws.Range("M132").Value = -2972.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("H132").Value = 7105.4443
$ws.Range("I132").Value = 7856.6787
$ws.Range("J132").Value = 4476.125
$ws.Range("K132").Value = 23570.0361
$ws.Range("L132").Value = 13428.375
$ws.Range("M132").Value = -21040.0361
$ws.Range("N132").Value = -18488.375
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0
